$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(69, 8).Value = 4678.6875
$ws.Cells.Item(69, 9).Value = 4667
$ws.Cells.Item(69, 10).Value = 4685.7
$ws.Cells.Item(69, 11).Value = 14001
$ws.Cells.Item(69, 12).Value = 14057.1
$ws.Cells.Item(69, 13).Value = -13127
$ws.Cells.Item(69, 14).Value = -15805.1
$ws.Cells.Item(70, 8).Value = 1671.8667
$ws.Cells.Item(70, 10).Value = 1591.2858
$ws.Cells.Item(70, 12).Value = 4773.857400000001
$ws.Cells.Item(70, 14).Value = -5313.857400000001
$ws.Cells.Item(72, 8).Value = 4678.6875
$ws.Cells.Item(72, 9).Value = 4667
$ws.Cells.Item(72, 10).Value = 4685.7
$ws.Cells.Item(72, 11).Value = 42003
$ws.Cells.Item(72, 12).Value = 42171.3
$ws.Cells.Item(72, 13).Value = -37635
$ws.Cells.Item(72, 14).Value = -50907.3
$ws.Cells.Item(73, 8).Value = 1671.8667
$ws.Cells.Item(73, 10).Value = 1591.2858
$ws.Cells.Item(73, 12).Value = 4773.857400000001
$ws.Cells.Item(73, 14).Value = -6645.857400000001
$ws.Cells.Item(112, 8).Value = 1660.4814
$ws.Cells.Item(112, 10).Value = 1863.6522
$ws.Cells.Item(112, 12).Value = 5590.9566
$ws.Cells.Item(112, 14).Value = -7806.9566
$ws.Cells.Item(137, 8).Value = 2083.103
$ws.Cells.Item(137, 9).Value = 1203.2954
$ws.Cells.Item(137, 10).Value = 3696.0833
$ws.Cells.Item(137, 11).Value = 3609.8862
$ws.Cells.Item(137, 12).Value = 11088.2499
$ws.Cells.Item(137, 13).Value = -1059.8862
$ws.Cells.Item(137, 14).Value = -16188.2499

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 18305.174
$ws.Cells.Item(32, 9).Value = 19150.408
$ws.Cells.Item(32, 10).Value = 12106.777
$ws.Cells.Item(32, 11).Value = 19150.408
$ws.Cells.Item(32, 12).Value = 12106.777
$ws.Cells.Item(32, 13).Value = -18863.408
$ws.Cells.Item(32, 14).Value = -12680.777
$ws.Cells.Item(102, 8).Value = 2241.5
$ws.Cells.Item(102, 9).Value = 1487.0714
$ws.Cells.Item(102, 10).Value = 4001.8333
$ws.Cells.Item(102, 11).Value = 1487.0714
$ws.Cells.Item(102, 12).Value = 4001.8333
$ws.Cells.Item(102, 13).Value = 134.9286
$ws.Cells.Item(102, 14).Value = -7245.8333
$ws.Cells.Item(128, 8).Value = 62641.25
$ws.Cells.Item(128, 10).Value = 62641.25
$ws.Cells.Item(128, 12).Value = 62641.25
$ws.Cells.Item(128, 14).Value = -72601.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(132, 8).Value = 47186.766
$ws.Cells.Item(132, 9).Value = 20000
$ws.Cells.Item(132, 10).Value = 55551.92
$ws.Cells.Item(132, 11).Value = 20000
$ws.Cells.Item(132, 12).Value = 55551.92
$ws.Cells.Item(132, 13).Value = -14940
$ws.Cells.Item(132, 14).Value = -65671.92
$ws.Cells.Item(134, 8).Value = 34216.938
$ws.Cells.Item(134, 9).Value = 3150.7368
$ws.Cells.Item(134, 10).Value = 79621.38
$ws.Cells.Item(134, 11).Value = 9452.2104
$ws.Cells.Item(134, 12).Value = 238864.14
$ws.Cells.Item(134, 13).Value = -6917.2104
$ws.Cells.Item(134, 14).Value = -243934.14

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 252.125
$ws.Cells.Item(22, 9).Value = 204.75
$ws.Cells.Item(22, 11).Value = 204.75
$ws.Cells.Item(22, 13).Value = 145.25
$ws.Cells.Item(58, 8).Value = 1717574.9
$ws.Cells.Item(58, 9).Value = 2598929.2
$ws.Cells.Item(58, 10).Value = 3830
$ws.Cells.Item(58, 11).Value = 2598929.2
$ws.Cells.Item(58, 12).Value = 3830
$ws.Cells.Item(58, 13).Value = -2598726.2
$ws.Cells.Item(58, 14).Value = -4236
$ws.Cells.Item(99, 8).Value = 1926.2858
$ws.Cells.Item(99, 9).Value = 1573.25
$ws.Cells.Item(99, 11).Value = 1573.25
$ws.Cells.Item(99, 13).Value = -75.25
$ws.Cells.Item(116, 8).Value = 48000
$ws.Cells.Item(116, 10).Value = 48000
$ws.Cells.Item(116, 12).Value = 48000
$ws.Cells.Item(116, 14).Value = -57178
$ws.Cells.Item(126, 8).Value = 1926.2858
$ws.Cells.Item(126, 9).Value = 1573.25
$ws.Cells.Item(126, 11).Value = 4719.75
$ws.Cells.Item(126, 13).Value = -2249.75
$ws.Cells.Item(132, 8).Value = 2259.8086
$ws.Cells.Item(132, 9).Value = 1615.3462
$ws.Cells.Item(132, 10).Value = 3057.7144
$ws.Cells.Item(132, 11).Value = 4846.0386
$ws.Cells.Item(132, 12).Value = 9173.143199999999
$ws.Cells.Item(132, 13).Value = -2316.0386
$ws.Cells.Item(132, 14).Value = -14233.1432
$ws.Cells.Item(134, 8).Value = 2663.6611
$ws.Cells.Item(134, 9).Value = 1630.2162
$ws.Cells.Item(134, 10).Value = 4401.727
$ws.Cells.Item(134, 11).Value = 4890.6486
$ws.Cells.Item(134, 12).Value = 13205.181
$ws.Cells.Item(134, 13).Value = -2355.6486
$ws.Cells.Item(134, 14).Value = -18275.181
$ws.Cells.Item(136, 8).Value = 1717574.9
$ws.Cells.Item(136, 9).Value = 2598929.2
$ws.Cells.Item(136, 10).Value = 3830
$ws.Cells.Item(136, 11).Value = 7796787.600000001
$ws.Cells.Item(136, 12).Value = 11490
$ws.Cells.Item(136, 13).Value = -7794237.600000001
$ws.Cells.Item(136, 14).Value = -16590

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(92, 8).Value = 988.8889
$ws.Cells.Item(92, 10).Value = 1018
$ws.Cells.Item(92, 12).Value = 3054
$ws.Cells.Item(92, 14).Value = -5550
$ws.Cells.Item(137, 8).Value = 17960.918
$ws.Cells.Item(137, 9).Value = 1996.6666
$ws.Cells.Item(137, 10).Value = 25623.76
$ws.Cells.Item(137, 11).Value = 5989.9998
$ws.Cells.Item(137, 12).Value = 76871.28
$ws.Cells.Item(137, 13).Value = -889.9997999999996
$ws.Cells.Item(137, 14).Value = -87071.28
$ws.Cells.Item(140, 8).Value = 2737.879
$ws.Cells.Item(140, 9).Value = 2214.348
$ws.Cells.Item(140, 10).Value = 3942
$ws.Cells.Item(140, 11).Value = 6643.044
$ws.Cells.Item(140, 12).Value = 11826
$ws.Cells.Item(140, 13).Value = -1463.044
$ws.Cells.Item(140, 14).Value = -22186

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 5119.75
$ws.Cells.Item(122, 9).Value = 5219.4287
$ws.Cells.Item(122, 11).Value = 15658.2861
$ws.Cells.Item(122, 13).Value = -13208.2861
$ws.Cells.Item(126, 8).Value = 2792.1667
$ws.Cells.Item(126, 9).Value = 1751
$ws.Cells.Item(126, 10).Value = 3833.3333
$ws.Cells.Item(126, 11).Value = 5253
$ws.Cells.Item(126, 12).Value = 11499.9999
$ws.Cells.Item(126, 13).Value = -2783
$ws.Cells.Item(126, 14).Value = -16439.9999
$ws.Cells.Item(132, 8).Value = 5585.0376
$ws.Cells.Item(132, 9).Value = 4269.41
$ws.Cells.Item(132, 10).Value = 9250
$ws.Cells.Item(132, 11).Value = 12808.23
$ws.Cells.Item(132, 12).Value = 27750
$ws.Cells.Item(132, 13).Value = -10278.23
$ws.Cells.Item(132, 14).Value = -32810

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1955.0625
$ws.Cells.Item(16, 9).Value = 1369.9
$ws.Cells.Item(16, 11).Value = 1369.9
$ws.Cells.Item(16, 13).Value = -1199.9
$ws.Cells.Item(82, 8).Value = 1406.4166
$ws.Cells.Item(82, 9).Value = 1150
$ws.Cells.Item(82, 10).Value = 1662.8334
$ws.Cells.Item(82, 11).Value = 1150
$ws.Cells.Item(82, 12).Value = 1662.8334
$ws.Cells.Item(82, 13).Value = -789
$ws.Cells.Item(82, 14).Value = -2384.8334
$ws.Cells.Item(85, 8).Value = 1406.4166
$ws.Cells.Item(85, 9).Value = 1150
$ws.Cells.Item(85, 10).Value = 1662.8334
$ws.Cells.Item(85, 11).Value = 1150
$ws.Cells.Item(85, 12).Value = 1662.8334
$ws.Cells.Item(85, 13).Value = 98
$ws.Cells.Item(85, 14).Value = -4158.8334
$ws.Cells.Item(132, 8).Value = 2786.6487
$ws.Cells.Item(132, 9).Value = 3204.3958
$ws.Cells.Item(132, 10).Value = 2015.4231
$ws.Cells.Item(132, 11).Value = 9613.187399999999
$ws.Cells.Item(132, 12).Value = 6046.2693
$ws.Cells.Item(132, 13).Value = -7083.187399999999
$ws.Cells.Item(132, 14).Value = -11106.2693

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 3478.5715
$ws.Cells.Item(62, 9).Value = 3400
$ws.Cells.Item(62, 10).Value = 3517.8572
$ws.Cells.Item(62, 11).Value = 3400
$ws.Cells.Item(62, 12).Value = 3517.8572
$ws.Cells.Item(62, 13).Value = -2776
$ws.Cells.Item(62, 14).Value = -4765.8572
$ws.Cells.Item(65, 8).Value = 3478.5715
$ws.Cells.Item(65, 9).Value = 3400
$ws.Cells.Item(65, 10).Value = 3517.8572
$ws.Cells.Item(65, 11).Value = 17000
$ws.Cells.Item(65, 12).Value = 17589.286
$ws.Cells.Item(65, 13).Value = -13880
$ws.Cells.Item(65, 14).Value = -23829.286
$ws.Cells.Item(132, 8).Value = 1399.5636
$ws.Cells.Item(132, 9).Value = 562.7646999999999
$ws.Cells.Item(132, 10).Value = 2754.3809
$ws.Cells.Item(132, 11).Value = 1688.2941
$ws.Cells.Item(132, 12).Value = 8263.1427
$ws.Cells.Item(132, 13).Value = 841.7059000000002
$ws.Cells.Item(132, 14).Value = -13323.1427
$ws.Cells.Item(136, 8).Value = 3324.8433
$ws.Cells.Item(136, 9).Value = 1616.4286
$ws.Cells.Item(136, 10).Value = 6868.222
$ws.Cells.Item(136, 11).Value = 4849.2858
$ws.Cells.Item(136, 12).Value = 20604.666
$ws.Cells.Item(136, 13).Value = -2299.2858
$ws.Cells.Item(136, 14).Value = -25704.666
